$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap order: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
$cellsA = @("G3","G6","G10","G11","G12","G13","G14","G15","G17","G18","G19","G20","G21","G22","G24","G26","G29","G32","G36","G37","G38","G39","G40","G41","G43","G44","G45","G46","G47","G48","G50","G52","G55","G58","G62","G63","G64","G65","G66","G67","G69","G70","G71","G72","G73","G74","G76","G78","G83","G84","G85","G86","G90","G92","G93","G94","G96","G99","G101","G109","G110","G111","G112","G116","G118","G119","G120","G122","G125","G127","G135","G136","G137","G138","G142","G144","G145","G146","G148","G151","G153")
foreach ($cell in $cellsA) {
    $ws.Range($cell).Value = "System, dnasr281@gmail.com"
}

# Swap order: "System, system, backup@backdoor.com" -> "system, System, backup@backdoor.com"
$cellsB = @("G2","G28","G54")
foreach ($cell in $cellsB) {
    $ws.Range($cell).Value = "system, System, backup@backdoor.com"
}
